$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, matching the style of the existing header row (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Per-row time_taken values (plain, unstyled, matching data rows)
$ws.Range("F2").Value = "2021-10-05 10:51:30.481370"
$ws.Range("F3").Value = "2021-10-05 10:51:30.481380"
$ws.Range("F4").Value = "2021-10-05 10:51:30.481384"
$ws.Range("F5").Value = "2021-10-05 10:51:30.481386"
$ws.Range("F6").Value = "2021-10-05 10:51:30.481389"
$ws.Range("F7").Value = "2021-10-05 10:51:30.481392"
$ws.Range("F8").Value = "2021-10-05 10:51:30.481394"
$ws.Range("F9").Value = "2021-10-05 10:51:30.481397"
$ws.Range("F10").Value = "2021-10-05 10:51:30.481399"
$ws.Range("F11").Value = "2021-10-05 10:51:30.481402"
$ws.Range("F12").Value = "2021-10-05 10:51:30.481404"
$ws.Range("F13").Value = "2021-10-05 10:51:30.481407"
$ws.Range("F14").Value = "2021-10-05 10:51:30.481409"
$ws.Range("F15").Value = "2021-10-05 10:51:30.481412"
$ws.Range("F16").Value = "2021-10-05 10:51:30.481414"
$ws.Range("F17").Value = "2021-10-05 10:51:30.481417"
$ws.Range("F18").Value = "2021-10-05 10:51:30.481420"
$ws.Range("F19").Value = "2021-10-05 10:51:30.481422"
$ws.Range("F20").Value = "2021-10-05 10:51:30.481425"
$ws.Range("F21").Value = "2021-10-05 10:51:30.481427"
$ws.Range("F22").Value = "2021-10-05 10:51:30.481430"
$ws.Range("F23").Value = "2021-10-05 10:51:30.481432"
$ws.Range("F24").Value = "2021-10-05 10:51:30.481435"
$ws.Range("F25").Value = "2021-10-05 10:51:30.481437"
$ws.Range("F26").Value = "2021-10-05 10:51:30.481440"
$ws.Range("F27").Value = "2021-10-05 10:51:30.481442"
$ws.Range("F28").Value = "2021-10-05 10:51:30.481445"
$ws.Range("F29").Value = "2021-10-05 10:51:30.481447"
$ws.Range("F30").Value = "2021-10-05 10:51:30.481450"
$ws.Range("F31").Value = "2021-10-05 10:51:30.481452"
$ws.Range("F32").Value = "2021-10-05 10:51:30.481455"
$ws.Range("F33").Value = "2021-10-05 10:51:30.481457"
$ws.Range("F34").Value = "2021-10-05 10:51:30.481460"
$ws.Range("F35").Value = "2021-10-05 10:51:30.481463"
$ws.Range("F36").Value = "2021-10-05 10:51:30.481465"
$ws.Range("F37").Value = "2021-10-05 10:51:30.481467"
$ws.Range("F38").Value = "2021-10-05 10:51:30.481470"
$ws.Range("F39").Value = "2021-10-05 10:51:30.481472"
$ws.Range("F40").Value = "2021-10-05 10:51:30.481475"
$ws.Range("F41").Value = "2021-10-05 10:51:30.481477"
